$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.466.24'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '3.019.00'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'" + '596.48'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').Value = "'" + '150.32'
$ws.Range('E6').Value = '  +6.05%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.016.11'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +10.86%  '
$ws.Range('E11').Value = '  +4.66%  '
$ws.Range('E13').Value = '  +3.24%  '
$ws.Range('D14').Value = "'" + '34.59'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D16').Value = '3.519.52'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').Value = '62.477.41'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('D18').Value = "'" + '7.01'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').Value = '3.020.07'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = "'" + '449.14'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('D23').Value = "'" + '7.45'
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').Value = "'" + '10.92'
$ws.Range('E25').Value = '  +11.37%  '
$ws.Range('E26').Value = '  +4.36%  '
$ws.Range('D27').Value = "'" + '12.08'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  +3.18%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = "'" + '7.30'
$ws.Range('E30').Value = '  +6.74%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').Value = "'" + '1.00'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('E32').Value = '  +4.20%  '
$ws.Range('D33').Value = "'" + '27.52'
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  +10.03%  '
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('E38').Value = '  +8.77%  '
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = "'" + '50.09'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = "'" + '9.05'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('E44').Value = '  +8.82%  '
$ws.Range('D45').Value = "'" + '390.59'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('D46').Value = "'" + '0.0354'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '2.736.52'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('D48').Value = "'" + '132.73'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = "'" + '2.18'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('E51').Value = '  +0.03%  '
